$d = $word.ActiveDocument

function ReplaceAllText($old, $new) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    [void]$range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- 1. Merge split runs that represent dates (no textual change, just run-merge) ---
ReplaceAllText "20/11/18" "20/11/18"
ReplaceAllText "21/11/18" "21/11/18"
ReplaceAllText "22/11/18" "22/11/18"
ReplaceAllText "23/11/18" "23/11/18"
ReplaceAllText "24/11/18" "24/11/18"
ReplaceAllText "25/11/18" "25/11/18"
ReplaceAllText "26/11/18" "26/11/18"
ReplaceAllText "27/11/18" "27/11/18"

# --- 2. Merge split runs for task-name cells (no textual change, just run-merge) ---
ReplaceAllText "Especificación de CU Generar PDF Programa" "Especificación de CU Generar PDF Programa"
ReplaceAllText "Codificación de CU Visualizar Programa" "Codificación de CU Visualizar Programa"
ReplaceAllText "Especificación de CU Filtrado y Visualización de Programas (móvil) " "Especificación de CU Filtrado y Visualización de Programas (móvil) "
ReplaceAllText "Diagramas del CU Filtrado y Visualización de Programas (móvil)" "Diagramas del CU Filtrado y Visualización de Programas (móvil)"

# --- 3. Locate the empty "PSI - Título 1" paragraph right before the "Recursos"
#        heading (it sits right after the "27/11 Fin de la segunda..." paragraph),
#        change its style to "PSI - Normal", and insert the new section
#        (closing planning table + two spacer paragraphs) right after it. ---

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.Trim() -eq "" -and $para.Style.NameLocal -eq "PSI - T$([char]0x00ED)tulo 1") {
        $next = $d.Paragraphs($i + 1)
        if ($next.Range.Text.Trim() -eq "Recursos") {
            $anchorIndex = $i
            break
        }
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate anchor paragraph"
}

$anchorPara = $d.Paragraphs($anchorIndex)
$anchorPara.Style = "PSI - Normal"

$insertXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="PSI-Ttulo2"/></w:pPr><w:r><w:t>Planificación etapa final de la iteración – fin de cuatrimestre:</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="Tablaconcuadrcula"/><w:tblW w:w="9820" w:type="dxa"/><w:jc w:val="center"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="01E0" w:firstRow="1" w:lastRow="1" w:firstColumn="1" w:lastColumn="1" w:noHBand="0" w:noVBand="0"/></w:tblPr><w:tblGrid><w:gridCol w:w="6980"/><w:gridCol w:w="2840"/></w:tblGrid><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="E6E6E6"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Nombre de la Tarea</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="E6E6E6"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Responsable</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Revisión de especificación, diagrama e implementación de CU Generar PDF</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Francisco Estrada</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Revisión de diagramas en CU previamente implementados</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Francisco Estrada</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Gestión y Análisis de Riesgos</w:t></w:r><w:r><w:t xml:space="preserve"> (incluyendo riesgos a futuro)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Francisco Estrada</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t xml:space="preserve">Finalización del CU Visualizar </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Programa  PDF</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Francisco Estrada</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Plan de Pruebas</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Nicolás Sartini</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Arquitectura del sistema</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Nicolás Sartini</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t xml:space="preserve">Revisión de especificación, diagrama e implementación de </w:t></w:r><w:r><w:t>aplicación móvil</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Fabricio González</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Revisión de GUI de aplicación móvil</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Fabricio González</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Planificación a largo plazo para finalización de proyecto</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Fabricio González</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Armado de presentación para martes 27/11</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Fabricio González</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Estimación</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Fabricio González</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Revisión de Diagrama de CU (incorporar CU de bibliografía)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Nicolás Sartini</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Balance de cursada y opinión personal del equipo</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t xml:space="preserve">Todo el </w:t></w:r><w:bookmarkStart w:id="6" w:name="_GoBack"/><w:bookmarkEnd w:id="6"/><w:r><w:t>equipo</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Grabar videos de los CU implementados en MP4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Todo el equipo</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4395" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Revi</w:t></w:r><w:r><w:t>sión</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>de documentación faltante</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1788" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="PSI-Normal"/></w:pPr><w:r><w:t>Nicolás Sartini</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:pStyle w:val="PSI-Ttulo2"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PSI-Ttulo1"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange = $d.Paragraphs($anchorIndex).Range
$insertRange.Collapse(0)
[void]$insertRange.InsertXML($insertXml)

Write-Output "structural insert done"
